# Extend the "appeals to the Ombudsman" table with three more years of
# data (2021, 2022, 2023) in new columns R, S, T, mirroring the existing
# per-row formatting used by column Q (the previous last year, 2020).
#
# Row 2 is a blank (bottom-bordered) spacer row -> new cells stay empty.
# Row 3 holds the year headers.
# Row 4 holds "written appeals" counts.
# Row 5 holds "positively resolved" counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column Q's cell formatting/styles down into the three new columns
# for rows 2-5 (this also copies Q's values, which are overwritten below).
$ws.Range("Q2:Q5").Copy($ws.Range("R2:R5"))
$ws.Range("Q2:Q5").Copy($ws.Range("S2:S5"))
$ws.Range("Q2:Q5").Copy($ws.Range("T2:T5"))

# Row 3: year headers
$ws.Range("R3").Value = 2021
$ws.Range("S3").Value = 2022
$ws.Range("T3").Value = 2023

# Row 4: number of written appeals per year
$ws.Range("R4").Value = 4301
$ws.Range("S4").Value = 3690
$ws.Range("T4").Value = 2620

# Row 5: number of positively resolved appeals per year
$ws.Range("R5").Value = 427
$ws.Range("S5").Value = 280
$ws.Range("T5").Value = 264

# Row 2 (spacer row) keeps the copied styling with no values, matching R2/S2/T2.

# Reset the active selection back to the sheet's home cell.
$ws.Range("A1").Select()
